$d = $word.ActiveDocument

# 1. "The version used on this heuristic substracts 3 to X and Y, what give us
#    sqrt( (x-3)^2 + (y-3)^2) since the center if the grid is at 3x3"
#    -> "The version used on this heuristic subtracts the half of the board
#    size (h) to X and Y, what give us sqrt( (x-h)^2 + (y-h)^2) since the
#    center if the grid is at 3x3"
$r1 = $d.Content.Find.Execute(
    "The version used on this heuristic substracts 3 to X and Y, what give us sqrt( (x-3)^2 + (y-3)^2) since the center if the grid is at 3x3",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The version used on this heuristic subtracts the half of the board size (h) to X and Y, what give us sqrt( (x-h)^2 + (y-h)^2) since the center if the grid is at 3x3",
    2)
Write-Output "edit1: $r1"

# 2. "moves available for the players with respect to all available moves"
#    -> "moves available for each player with respect to all available moves"
$r2 = $d.Content.Find.Execute(
    "moves available for the players with respect to all available moves",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "moves available for each player with respect to all available moves",
    2)
Write-Output "edit2: $r2"

# 3. Add a new sentence right after the Proportion Score heuristic's formula
#    paragraph ("my_proportion * 10 - opponent_proportion * 10"): the
#    paragraph directly following it is blank, and the new sentence becomes
#    its own paragraph right after that blank one (before the next blank
#    paragraph that precedes "Results:").
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "my_proportion \* 10 - opponent_proportion \* 10") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    $blank = $d.Paragraphs.Item($target + 1)
    $blank.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($target + 2)
    $newPara.Range.Text = "This heuristic worked good on average, however I got some looses more frequently than H1"
    Write-Output "edit3: inserted after paragraph $target"
} else {
    Write-Output "edit3: target paragraph not found"
}
